# Applies the edit described by the commit:
#   - Cell B6 on Sheet1 changes from "240.63.255.255." to "sdaasf"
#     (the old shared string becomes unused and is dropped; the new
#     string "sdaasf" is appended to the shared-string table, which is
#     exactly what Excel itself does when the last use of a shared
#     string is replaced with new text).
#   - The active selection / cursor moves to cell B6 (and the previous
#     "topLeftCell" scroll-position override is cleared as a natural
#     side effect of re-selecting the cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "sdaasf"
[void]$ws.Range("B6").Select()
